# Generate Report for Archive
#
# The localization-status report was regenerated. For the entry
# "a80a1b69-ced6-4051-8b8f-5300c3d93bfd.md" the Status moved from
# "Ready for handoff" to "In Translation". As a side effect, the row
# order of "a80a1b69-...md" and "fe8086df-...md" swapped (rows 6 and 7)
# on every sheet, while all other rows (and the hyperlink targets behind
# each rId) stayed exactly where they were.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

# Row 6 becomes the a80a1b69 entry, row 7 becomes the fe8086df entry.
$ws.Range("A6").Value2 = "a80a1b69-ced6-4051-8b8f-5300c3d93bfd.md"
$ws.Range("B6").Value2 = "e2e\a80a1b69-ced6-4051-8b8f-5300c3d93bfd.md"
$ws.Range("E6").Value2 = "In Translation"
$ws.Range("F6").Value2 = "In Translation"
$ws.Range("G6").Value2 = "2016-09-05 22:55:37"

$ws.Range("A7").Value2 = "fe8086df-5963-4f5a-a467-1b24aa2d5d34.md"
$ws.Range("B7").Value2 = "e2e\fe8086df-5963-4f5a-a467-1b24aa2d5d34.md"
$ws.Range("E7").Value2 = "In Translation"
$ws.Range("F7").Value2 = "In Translation"
$ws.Range("G7").Value2 = "2016-09-05 22:52:50"

foreach ($h in $ws.Hyperlinks) {
    if ($h.Range.Row -eq 6 -and $h.Range.Column -eq 2) {
        $h.TextToDisplay = "e2e\a80a1b69-ced6-4051-8b8f-5300c3d93bfd.md"
    }
    elseif ($h.Range.Row -eq 7 -and $h.Range.Column -eq 2) {
        $h.TextToDisplay = "e2e\fe8086df-5963-4f5a-a467-1b24aa2d5d34.md"
    }
}

# ---------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A6").Value2 = "a80a1b69-ced6-4051-8b8f-5300c3d93bfd.md"
$ws.Range("C6").Value2 = "In Translation"
$ws.Range("G6").Value2 = "a80a1b69-ced6-4051-8b8f-5300c3d93bfd.a81fe4cea0ccef8e6d177f145c89ada16bb5108b.zh-cn.xlf"
$ws.Range("H6").Value2 = "2016-09-05 22:55:32"

$ws.Range("A7").Value2 = "fe8086df-5963-4f5a-a467-1b24aa2d5d34.md"
$ws.Range("C7").Value2 = "In Translation"
$ws.Range("G7").Value2 = "fe8086df-5963-4f5a-a467-1b24aa2d5d34.4ac973d788bf32ce935fb0c4ffdae4c4b8bb96aa.zh-cn.xlf"
$ws.Range("H7").Value2 = "2016-09-05 22:52:46"

foreach ($h in $ws.Hyperlinks) {
    if ($h.Range.Row -eq 6 -and $h.Range.Column -eq 1) {
        $h.TextToDisplay = "a80a1b69-ced6-4051-8b8f-5300c3d93bfd.md"
    }
    elseif ($h.Range.Row -eq 7 -and $h.Range.Column -eq 1) {
        $h.TextToDisplay = "fe8086df-5963-4f5a-a467-1b24aa2d5d34.md"
    }
}

# ---------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A6").Value2 = "a80a1b69-ced6-4051-8b8f-5300c3d93bfd.md"
$ws.Range("C6").Value2 = "In Translation"
$ws.Range("G6").Value2 = "a80a1b69-ced6-4051-8b8f-5300c3d93bfd.a81fe4cea0ccef8e6d177f145c89ada16bb5108b.de-de.xlf"
$ws.Range("H6").Value2 = "2016-09-05 22:55:37"

$ws.Range("A7").Value2 = "fe8086df-5963-4f5a-a467-1b24aa2d5d34.md"
$ws.Range("C7").Value2 = "In Translation"
$ws.Range("G7").Value2 = "fe8086df-5963-4f5a-a467-1b24aa2d5d34.4ac973d788bf32ce935fb0c4ffdae4c4b8bb96aa.de-de.xlf"
$ws.Range("H7").Value2 = "2016-09-05 22:52:50"

foreach ($h in $ws.Hyperlinks) {
    if ($h.Range.Row -eq 6 -and $h.Range.Column -eq 1) {
        $h.TextToDisplay = "a80a1b69-ced6-4051-8b8f-5300c3d93bfd.md"
    }
    elseif ($h.Range.Row -eq 7 -and $h.Range.Column -eq 1) {
        $h.TextToDisplay = "fe8086df-5963-4f5a-a467-1b24aa2d5d34.md"
    }
}
